$wb = $excel.ActiveWorkbook
$win = $wb.Windows.Item(1)
Write-Output $win.WindowState
$win.WindowState = -4137
